$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store a literal text value, even when the
# string looks like a number (matches the source inlineStr text cells).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "69.105.35"
Set-TextValue "D3" "3.811.01"
Set-TextValue "E3" "  +1.84%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "601.16"
Set-TextValue "E5" "  -0.01%  "
Set-TextValue "D6" "164.22"
Set-TextValue "E6" "  -2.24%  "
Set-TextValue "D7" "3.810.41"
Set-TextValue "E7" "  +1.89%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "E9" "  -0.47%  "
Set-TextValue "E10" "  +2.05%  "
Set-TextValue "D11" "6.34"
Set-TextValue "E11" "  -1.55%  "
Set-TextValue "E12" "  +0.07%  "
Set-TextValue "D13" "36.97"
Set-TextValue "E13" "  -2.60%  "
Set-TextValue "D14" "0.0000245"
Set-TextValue "E14" "  -0.62%  "
Set-TextValue "D15" "4.444.12"
Set-TextValue "E15" "  +1.90%  "
Set-TextValue "D16" "3.835.37"
Set-TextValue "E16" "  +2.66%  "
Set-TextValue "D17" "69.238.71"
Set-TextValue "E17" "  +0.23%  "
Set-TextValue "D18" "7.52"
Set-TextValue "E18" "  +3.22%  "
Set-TextValue "E19" "  +7.37%  "
Set-TextValue "E20" "  +0.11%  "
Set-TextValue "D21" "17.22"
Set-TextValue "E21" "  +0.57%  "
Set-TextValue "D22" "487.34"
Set-TextValue "E22" "  -0.95%  "
Set-TextValue "D23" "0.720"
Set-TextValue "E23" "  -0.50%  "
Set-TextValue "D24" "0.0000156"
Set-TextValue "E24" "  +4.38%  "
Set-TextValue "D25" "84.43"
Set-TextValue "E25" "  -0.34%  "
Set-TextValue "E26" "  -2.32%  "
Set-TextValue "D27" "12.14"
Set-TextValue "E27" "  -1.21%  "
Set-TextValue "D28" "10.04"
Set-TextValue "E28" "  -0.74%  "
Set-TextValue "E29" "  +0.02%  "
Set-TextValue "D30" "2.97"
Set-TextValue "E30" "  -0.34%  "
Set-TextValue "E31" "  -0.28%  "
Set-TextValue "E32" "  -3.99%  "
Set-TextValue "D33" "3.963.02"
Set-TextValue "E33" "  +2.05%  "
Set-TextValue "D34" "31.87"
Set-TextValue "E34" "  +1.14%  "
Set-TextValue "D35" "3.750.61"
Set-TextValue "E35" "  +2.13%  "
Set-TextValue "E36" "  -1.79%  "
Set-TextValue "D37" "1.03"
Set-TextValue "E37" "  +0.87%  "
Set-TextValue "E38" "  +4.41%  "
Set-TextValue "E39" "  +0.36%  "
Set-TextValue "E40" "  +0.04%  "
Set-TextValue "B41" "TheGraph"
Set-TextValue "C41" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D41" "0.318"
Set-TextValue "E41" "  -1.48%  "
Set-TextValue "B42" "dogwifhat"
Set-TextValue "C42" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D42" "3.03"
Set-TextValue "E42" "  +2.82%  "
Set-TextValue "D43" "434.47"
Set-TextValue "E43" "  +0.99%  "
Set-TextValue "E44" "  +0.06%  "
Set-TextValue "E45" "  +0.39%  "
Set-TextValue "D47" "8.37"
Set-TextValue "E47" "  -1.13%  "
Set-TextValue "D48" "143.00"
Set-TextValue "E48" "  +1.04%  "
Set-TextValue "D49" "2.824.89"
Set-TextValue "E49" "  +1.73%  "
Set-TextValue "D50" "0.0354"
Set-TextValue "E50" "  +0.29%  "
Set-TextValue "D51" "39.23"
Set-TextValue "E51" "  -1.89%  "
